$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.188892841339111
$ws.Range("B1").Value = 2.194520950317383
$ws.Range("C1").Value = 3.485339641571045
$ws.Range("D1").Value = 2.105850219726562
$ws.Range("E1").Value = 1.082806825637817
